$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 4, 6) {
    # Force the F column to stay as text (rather than being reinterpreted
    # as a date serial number) before writing the literal string value.
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = "2024-08-04"

    $ws.Cells.Item($r, 9).Value = "Yes"
    $ws.Cells.Item($r, 10).Value = "Yes"
}
